$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of rows 2-3 with rows 4-5 (the "Golden Nugget / Especial / Primera"
# record pair moves from 2022-12-07 to 2021-11-05, and the "Californiana(o) / Primera" pair
# moves from 2021-11-05 to 2022-12-07), while columns A,B,C,E,F,G,H,I,J,R,T stay identical.

# --- Row 2 (was Golden Nugget/Especial @44902) becomes Californiana(o)/Primera @44505 ---
$ws.Range("D2").Value = 44505
$ws.Range("K2").Value = "Californiana(o)"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("S2").Value = 1500

# --- Row 3 (was Golden Nugget/Primera @44902) becomes Golden Nugget/Primera @44505 (qty/price swap) ---
$ws.Range("D3").Value = 44505
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1500

# --- Row 4 (was Californiana(o)/Primera @44505) becomes Golden Nugget/Especial @44902 ---
$ws.Range("D4").Value = 44902
$ws.Range("K4").Value = "Golden Nugget"
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/caja 10 kilos"
$ws.Range("S4").Value = 1500

# --- Row 5 (was Golden Nugget/Primera @44505) becomes Golden Nugget/Primera @44902 (qty/price swap) ---
$ws.Range("D5").Value = 44902
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("Q5").Value = "$/caja 10 kilos"
$ws.Range("S5").Value = 1300
